$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.968.24"
$ws.Range("E2").Value = "  +4.24%  "
$ws.Range("D3").Value = "2.644.84"
$ws.Range("E3").Value = "  +6.36%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'113.71"
$ws.Range("E5").Value = "  +8.53%  "
$ws.Range("D6").Value = "'326.64"
$ws.Range("E6").Value = "  +3.36%  "
$ws.Range("E7").Value = "  +2.34%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "'0.556"
$ws.Range("E9").Value = "  +4.52%  "
$ws.Range("D10").Value = "'41.09"
$ws.Range("E10").Value = "  +6.91%  "
$ws.Range("D11").Value = "'20.15"
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("E12").Value = "  +3.06%  "
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("D14").Value = "'7.39"
$ws.Range("E14").Value = "  +5.31%  "
$ws.Range("D15").Value = "3.061.42"
$ws.Range("E15").Value = "  +6.49%  "
$ws.Range("D16").Value = "2.640.55"
$ws.Range("E16").Value = "  +6.03%  "
$ws.Range("D17").Value = "'0.872"
$ws.Range("E17").Value = "  +6.00%  "
$ws.Range("D18").Value = "49.892.15"
$ws.Range("E18").Value = "  +4.56%  "
$ws.Range("D19").Value = "'13.18"
$ws.Range("E19").Value = "  +2.92%  "
$ws.Range("E20").Value = "  +3.68%  "
$ws.Range("D21").Value = "'2.91"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("D22").Value = "0.0₃0956"
$ws.Range("E22").Value = "  +3.46%  "
$ws.Range("E23").Value = "  +1.75%  "
$ws.Range("D24").Value = "'276.11"
$ws.Range("E24").Value = "  +2.78%  "
$ws.Range("E25").Value = "  +3.81%  "
$ws.Range("D26").Value = "'26.81"
$ws.Range("E26").Value = "  +4.56%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "'10.03"
$ws.Range("E28").Value = "  +4.10%  "
$ws.Range("E29").Value = "  +1.73%  "
$ws.Range("D30").Value = "'36.16"
$ws.Range("E30").Value = "  +5.51%  "
$ws.Range("D31").Value = "'0.142"
$ws.Range("E31").Value = "  +2.85%  "
$ws.Range("D32").Value = "'50.26"
$ws.Range("E32").Value = "  +2.18%  "
$ws.Range("D33").Value = "'5.48"
$ws.Range("E33").Value = "  +4.20%  "
$ws.Range("D34").Value = "'19.45"
$ws.Range("E34").Value = "  +3.03%  "
$ws.Range("E35").Value = "  +5.19%  "
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("D37").Value = "'5.03"
$ws.Range("E37").Value = "  +10.75%  "
$ws.Range("E38").Value = "  +7.53%  "
$ws.Range("D39").Value = "'3.10"
$ws.Range("E39").Value = "  +8.65%  "
$ws.Range("D40").Value = "'123.73"
$ws.Range("E40").Value = "  +1.46%  "
$ws.Range("E41").Value = "  +2.38%  "
$ws.Range("E42").Value = "  +0.64%  "
$ws.Range("D43").Value = "'22.03"
$ws.Range("E43").Value = "  -0.81%  "
$ws.Range("E44").Value = "  +5.19%  "
$ws.Range("D45").Value = "2.083.97"
$ws.Range("E45").Value = "  +4.69%  "
$ws.Range("E46").Value = "  +7.17%  "
$ws.Range("D47").Value = "'2.31"
$ws.Range("E47").Value = "  +15.90%  "
$ws.Range("E48").Value = "  +5.92%  "
$ws.Range("D49").Value = "'9.12"
$ws.Range("E49").Value = "  +3.09%  "
$ws.Range("E50").Value = "  +5.74%  "
$ws.Range("D51").Value = "'59.60"
$ws.Range("E51").Value = "  +6.53%  "
